# "note algoritmi + fisiere intrebari examen"
# Fill in the "sapt. 14" (P) attendance column with 1 for each student that
# still needed it, and record the exam grade ("Nota", column R) for every
# student who already has a full attendance record. Column Q ("Prezente")
# is a SUM(C:P) formula, so it recalculates automatically once P is set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row -> (P14 attendance, R grade)
$updates = @{
    3  = @(1, 6)
    4  = @(1, 6)
    6  = @(1, 8)
    8  = @(1, 7)
    10 = @(1, 7)
    11 = @(1, 10)
    12 = @(1, 7)
    13 = @(1, 9)
    14 = @(1, 7)
    15 = @(1, 7)
    17 = @(1, 6)
    18 = @(1, 8)
    19 = @(1, 6)
    20 = @(1, 10)
    22 = @(1, 9)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 16).Value = $vals[0]   # column P = sapt. 14
    $ws.Cells.Item($row, 18).Value = $vals[1]   # column R = Nota
}

$wb.Application.Calculate()

# Leave the view the way the author left it: the last touched cell (R20)
# selected.
$ws.Range("R20").Select()
